$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.090786576271057
$ws.Range("B1").Value = 1.956804990768433
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.096278429031372
$ws.Range("E1").Value = 1.13579535484314
